$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Value = "Departamento de Química / Centro de Ciências Exatas e Tecnológicas / Universidade Federal de Viçosa - Campus Viçosa"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Value = "Departamento de Química / Centro de Ciências Exatas e Tecnológicas / Universidade Federal de Viçosa - Campus Viçosa"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Value = "Departamento de Química / Centro de Ciências Exatas e Tecnológicas / Universidade Federal de Viçosa - Campus Viçosa"
$ws.Range("J2").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null
$ws.Range("K2").Value = "Departamento de Química / Centro de Ciências Exatas e Tecnológicas / Universidade Federal de Viçosa - Campus Viçosa"
# Row 3
$ws.Range("D3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Value = "Departamento de Química  / Instituto de Ciências Exatas  / Universidade Federal de Minas Gerais"
# Row 4
$ws.Range("D4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Value = "Centre national de la recherche scientifique (CNRS)"
$ws.Range("F4").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4122) | Out-Null
$ws.Range("G4").Value = "Empresa Brasileira de Pesquisa Agropecuária"
$ws.Range("H4").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null
$ws.Range("I4").Value = "Departamento de Química  / Instituto de Ciências Exatas  / Universidade Federal de Minas Gerais"
# Row 5
$ws.Range("D5").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Value = "Department of Physical and Environmental Sciences / Environmental NMR Centre / University of Toronto"
$ws.Range("F5").Copy() | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null
$ws.Range("G5").Value = "Environmental NMR Centre / University of Toronto"
$ws.Range("H5").Copy() | Out-Null
$ws.Range("I5").PasteSpecial(-4122) | Out-Null
$ws.Range("I5").Value = "Synex Medical"
$ws.Range("J5").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null
$ws.Range("K5").Value = "Chemistry / Scarborough Campus / University of Toronto"
$ws.Range("L5").Copy() | Out-Null
$ws.Range("M5").PasteSpecial(-4122) | Out-Null
$ws.Range("M5").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("N5").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null
$ws.Range("O5").Value = "Chemistry / Environmental NMR Centre / University of Toronto"
$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial(-4122) | Out-Null
$ws.Range("Q5").Value = "Chemistry / Scarborough Campus / University of Toronto"
$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null
$ws.Range("S5").Value = "Environmental NMR Centre / University of Toronto"
$ws.Range("T5").Copy() | Out-Null
$ws.Range("U5").PasteSpecial(-4122) | Out-Null
$ws.Range("U5").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("V5").Copy() | Out-Null
$ws.Range("W5").PasteSpecial(-4122) | Out-Null
$ws.Range("W5").Value = "Environmental NMR Centre / University of Toronto"
$ws.Range("X5").Copy() | Out-Null
$ws.Range("Y5").PasteSpecial(-4122) | Out-Null
$ws.Range("Y5").Value = "Environmental Monitoring & Reporting Branch / Ontario Ministry of the Environment"
$ws.Range("Z5").Copy() | Out-Null
$ws.Range("AA5").PasteSpecial(-4122) | Out-Null
$ws.Range("AA5").Value = "Bruker"
$ws.Range("AB5").Copy() | Out-Null
$ws.Range("AC5").PasteSpecial(-4122) | Out-Null
$ws.Range("AC5").Value = "Bruker"
$ws.Range("AD5").Copy() | Out-Null
$ws.Range("AE5").PasteSpecial(-4122) | Out-Null
$ws.Range("AE5").Value = "Bruker"
$ws.Range("AF5").Copy() | Out-Null
$ws.Range("AG5").PasteSpecial(-4122) | Out-Null
$ws.Range("AG5").Value = "Department of Physical and Environmental Sciences / Scarborough Campus / University of Toronto"
$ws.Range("AH5").Copy() | Out-Null
$ws.Range("AI5").PasteSpecial(-4122) | Out-Null
$ws.Range("AI5").Value = "University of Toronto"
# Row 6
$ws.Range("D6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Value = "Instituto de Química de São Carlos  / Universidade de São Paulo"
$ws.Range("F6").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
$ws.Range("G6").Value = "Instituto de Química de São Carlos  / Universidade de São Paulo"
$ws.Range("H6").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null
$ws.Range("I6").Value = "Instituto de Química de São Carlos/ Universidade de São Paulo"
$ws.Range("J6").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$ws.Range("K6").Value = "Empresa Brasileira de Pesquisa Agropecuária"
$ws.Range("L6").Copy() | Out-Null
$ws.Range("M6").PasteSpecial(-4122) | Out-Null
$ws.Range("M6").Value = "Instituto de Química de São Carlos  / Universidade de São Paulo"
# Row 7
$ws.Range("D7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Value = "Institute of Chemistry/Federal University of Goiás"
$ws.Range("F7").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").Value = "Institute of Chemistry/Federal University of Goiás"
$ws.Range("H7").Copy() | Out-Null
$ws.Range("I7").PasteSpecial(-4122) | Out-Null
$ws.Range("I7").Value = "Institute of Chemistry/Federal University of Goiás"
# Row 8
$ws.Range("D8").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = "Universidade Federal dos Vales do Jequitinhonha e Mucuri"
$ws.Range("F8").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null
$ws.Range("G8").Value = "Universidade Federal dos Vales do Jequitinhonha e Mucuri"
# Row 9
$ws.Range("D9").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").Value = "Chemistry / Scarborough Campus / University of Toronto"
$ws.Range("F9").Copy() | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null
$ws.Range("G9").Value = "Universidade Federal de São Carlos"
$ws.Range("H9").Copy() | Out-Null
$ws.Range("I9").PasteSpecial(-4122) | Out-Null
$ws.Range("I9").Value = "Biophysical Chemistry / Max Planck Institute / Max Planck Institute for Biophysical Chemistry"
$ws.Range("J9").Copy() | Out-Null
$ws.Range("K9").PasteSpecial(-4122) | Out-Null
$ws.Range("K9").Value = "Biophysical Chemistry / Max Planck Institute / Max Planck Institute for Biophysical Chemistry"
$ws.Range("L9").Copy() | Out-Null
$ws.Range("M9").PasteSpecial(-4122) | Out-Null
$ws.Range("M9").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("N9").Copy() | Out-Null
$ws.Range("O9").PasteSpecial(-4122) | Out-Null
$ws.Range("O9").Value = "University of Toronto"
# Row 10
$ws.Range("D10").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = "Universidade Estadual de Campinas"
$ws.Range("F10").Copy() | Out-Null
$ws.Range("G10").PasteSpecial(-4122) | Out-Null
$ws.Range("G10").Value = "Universidade Estadual de Campinas"
$ws.Range("H10").Copy() | Out-Null
$ws.Range("I10").PasteSpecial(-4122) | Out-Null
$ws.Range("I10").Value = "Universidade Estadual de Campinas"
# Row 11
$ws.Range("D11").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Value = "Universidade de São Paulo"
$ws.Range("F11").Copy() | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null
$ws.Range("G11").Value = "Universidade de São Paulo"
$ws.Range("H11").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null
$ws.Range("I11").Value = "Universidade de São Paulo"
$ws.Range("J11").Copy() | Out-Null
$ws.Range("K11").PasteSpecial(-4122) | Out-Null
$ws.Range("K11").Value = "Universidade de São Paulo"
$ws.Range("L11").Copy() | Out-Null
$ws.Range("M11").PasteSpecial(-4122) | Out-Null
$ws.Range("M11").Value = "Instituto de Física de São Carlos"
# Row 12
$ws.Range("D12").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = "Universidade Federal de Goiás"
$ws.Range("F12").Copy() | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null
$ws.Range("G12").Value = "Universidade Federal de Goiás"
$ws.Range("H12").Copy() | Out-Null
$ws.Range("I12").PasteSpecial(-4122) | Out-Null
$ws.Range("I12").Value = "Institute of Chemistry/Federal University of Goiás"
$ws.Range("J12").Copy() | Out-Null
$ws.Range("K12").PasteSpecial(-4122) | Out-Null
$ws.Range("K12").Value = "Instituto de Química / Universidade Federal de Goiás"
# Row 13
$ws.Range("D13").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").Value = "Universidade Federal do Rio Grande do Sul"
$ws.Range("F13").Copy() | Out-Null
$ws.Range("G13").PasteSpecial(-4122) | Out-Null
$ws.Range("G13").Value = "Universidade Federal de São Carlos"
$ws.Range("H13").Copy() | Out-Null
$ws.Range("I13").PasteSpecial(-4122) | Out-Null
$ws.Range("I13").Value = "Universidade Federal do Rio Grande do Sul"
# Row 14
$ws.Range("D14").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = "University of Toronto"
$ws.Range("F14").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("H14").Copy() | Out-Null
$ws.Range("I14").PasteSpecial(-4122) | Out-Null
$ws.Range("I14").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("J14").Copy() | Out-Null
$ws.Range("K14").PasteSpecial(-4122) | Out-Null
$ws.Range("K14").Value = "Institute of Microstructure Technology / Karlsruhe Institute of Technology / Karlsruhe Institute of Technology"
$ws.Range("L14").Copy() | Out-Null
$ws.Range("M14").PasteSpecial(-4122) | Out-Null
$ws.Range("M14").Value = "Institute of Microstructure Technology / Karlsruhe Institute of Technology / Karlsruhe Institute of Technology"
$ws.Range("N14").Copy() | Out-Null
$ws.Range("O14").PasteSpecial(-4122) | Out-Null
$ws.Range("O14").Value = "Institute of Microstructure Technology / Karlsruhe Institute of Technology / Karlsruhe Institute of Technology"
$ws.Range("P14").Copy() | Out-Null
$ws.Range("Q14").PasteSpecial(-4122) | Out-Null
$ws.Range("Q14").Value = "Institute of Microstructure Technology / Karlsruhe Institute of Technology / Karlsruhe Institute of Technology"
$ws.Range("R14").Copy() | Out-Null
$ws.Range("S14").PasteSpecial(-4122) | Out-Null
$ws.Range("S14").Value = "Institute of Microstructure Technology / Karlsruhe Institute of Technology / Karlsruhe Institute of Technology"
$ws.Range("T14").Copy() | Out-Null
$ws.Range("U14").PasteSpecial(-4122) | Out-Null
$ws.Range("U14").Value = "University of Toronto"
# Row 15
$ws.Range("D15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = "Universidade Estadual de Campinas"
$ws.Range("F15").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = "Universidade Estadual de Campinas"
$ws.Range("H15").Copy() | Out-Null
$ws.Range("I15").PasteSpecial(-4122) | Out-Null
$ws.Range("I15").Value = "Universidade Federal de Minas Gerais"
$ws.Range("J15").Copy() | Out-Null
$ws.Range("K15").PasteSpecial(-4122) | Out-Null
$ws.Range("K15").Value = "Universidade Federal de Minas Gerais"
$ws.Range("L15").Copy() | Out-Null
$ws.Range("M15").PasteSpecial(-4122) | Out-Null
$ws.Range("M15").Value = "Universidade Estadual de Campinas"
$ws.Range("N15").Copy() | Out-Null
$ws.Range("O15").PasteSpecial(-4122) | Out-Null
$ws.Range("O15").Value = "Universidade Estadual de Campinas"
$ws.Range("P15").Copy() | Out-Null
$ws.Range("Q15").PasteSpecial(-4122) | Out-Null
$ws.Range("Q15").Value = "Department of Chemistry / McGill University / McGill University"
$ws.Range("R15").Copy() | Out-Null
$ws.Range("S15").PasteSpecial(-4122) | Out-Null
$ws.Range("S15").Value = "Universidade Estadual de Campinas"
# Row 16
$ws.Range("D16").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = "Universidade Estadual de Campinas"
$ws.Range("F16").Copy() | Out-Null
$ws.Range("G16").PasteSpecial(-4122) | Out-Null
$ws.Range("G16").Value = "Universidade Estadual de Campinas"
# Row 17
$ws.Range("D17").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = "Universidade Federal de Ouro Preto"
$ws.Range("F17").Copy() | Out-Null
$ws.Range("G17").PasteSpecial(-4122) | Out-Null
$ws.Range("G17").Value = "Universidade Federal de Ouro Preto"
$ws.Range("H17").Copy() | Out-Null
$ws.Range("I17").PasteSpecial(-4122) | Out-Null
$ws.Range("I17").Value = "Universidade de São Paulo"
# Row 18
$ws.Range("D18").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = "Instituto Militar de Engenharia"
$ws.Range("F18").Copy() | Out-Null
$ws.Range("G18").PasteSpecial(-4122) | Out-Null
$ws.Range("G18").Value = "INCQS / Fundação Oswaldo Cruz"
$ws.Range("H18").Copy() | Out-Null
$ws.Range("I18").PasteSpecial(-4122) | Out-Null
$ws.Range("I18").Value = "INCQS/Fiocruz"
$ws.Range("J18").Copy() | Out-Null
$ws.Range("K18").PasteSpecial(-4122) | Out-Null
$ws.Range("K18").Value = "Instituto Militar de Engenharia"
$ws.Range("L18").Copy() | Out-Null
$ws.Range("M18").PasteSpecial(-4122) | Out-Null
$ws.Range("M18").Value = "Instituto Militar de Engenharia"
$ws.Range("N18").Copy() | Out-Null
$ws.Range("O18").PasteSpecial(-4122) | Out-Null
$ws.Range("O18").Value = "CDTS / Presidência / Fundação Oswaldo Cruz"
# Row 19
$ws.Range("D19").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").Value = "Universidade Estadual de Campinas"
$ws.Range("F19").Copy() | Out-Null
$ws.Range("G19").PasteSpecial(-4122) | Out-Null
$ws.Range("G19").Value = "Universidade Estadual de Campinas"
# Row 20
$ws.Range("D20").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "Federal Institute Goiano"
$ws.Range("F20").Copy() | Out-Null
$ws.Range("G20").PasteSpecial(-4122) | Out-Null
$ws.Range("G20").Value = "Institute of Chemistry, Federal University of Goias"
$ws.Range("H20").Copy() | Out-Null
$ws.Range("I20").PasteSpecial(-4122) | Out-Null
$ws.Range("I20").Value = "Federal Institute of Goiás"
$ws.Range("J20").Copy() | Out-Null
$ws.Range("K20").PasteSpecial(-4122) | Out-Null
$ws.Range("K20").Value = "Institute of Chemistry, Federal University of Goias"
# Row 21
$ws.Range("D21").Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").Value = "Universidade Federal de Ouro Preto"
$ws.Range("F21").Copy() | Out-Null
$ws.Range("G21").PasteSpecial(-4122) | Out-Null
$ws.Range("G21").Value = "Universidade Federal de Ouro Preto"
$ws.Range("H21").Copy() | Out-Null
$ws.Range("I21").PasteSpecial(-4122) | Out-Null
$ws.Range("I21").Value = "Universidade Federal de Ouro Preto"
$ws.Range("J21").Copy() | Out-Null
$ws.Range("K21").PasteSpecial(-4122) | Out-Null
$ws.Range("K21").Value = "Universidade Federal de Ouro Preto"
# Row 22
$ws.Range("D22").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "Instituto de Bioquímica Médica, Centro Nacional de Ressonância Magnética Nuclear Jiri Jonas, Universidade Federal do Rio de Janeiro"
$ws.Range("F22").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null
$ws.Range("G22").Value = "Departamento de Bioquímica, Instituto de Química, Universidade Federal do Rio de Janeiro"
$ws.Range("H22").Copy() | Out-Null
$ws.Range("I22").PasteSpecial(-4122) | Out-Null
$ws.Range("I22").Value = "Departamento de Bioquímica, Instituto de Química, Universidade Federal do Rio de Janeiro"
$ws.Range("J22").Copy() | Out-Null
$ws.Range("K22").PasteSpecial(-4122) | Out-Null
$ws.Range("K22").Value = "Departamento de Bioquímica, Instituto de Química, Universidade Federal do Rio de Janeiro"
$ws.Range("L22").Copy() | Out-Null
$ws.Range("M22").PasteSpecial(-4122) | Out-Null
$ws.Range("M22").Value = "Instituto de Bioquímica Médica, Centro Nacional de Ressonância Magnética Nuclear Jiri Jonas, Universidade Federal do Rio de Janeiro"
# Row 23
$ws.Range("D23").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "Universidade Federal do Rio de Janeiro"
$ws.Range("F23").Copy() | Out-Null
$ws.Range("G23").PasteSpecial(-4122) | Out-Null
$ws.Range("G23").Value = "Instituto de Bioquímica Médica, Centro Nacional de Ressonância Magnética Nuclear Jiri Jonas, Universidade Federal do Rio de Janeiro"
$ws.Range("H23").Copy() | Out-Null
$ws.Range("I23").PasteSpecial(-4122) | Out-Null
$ws.Range("I23").Value = "Universidade Federal do Rio de Janeiro"
$ws.Range("J23").Copy() | Out-Null
$ws.Range("K23").PasteSpecial(-4122) | Out-Null
$ws.Range("K23").Value = "Universidade Federal do Rio de Janeiro"
$ws.Range("L23").Copy() | Out-Null
$ws.Range("M23").PasteSpecial(-4122) | Out-Null
$ws.Range("M23").Value = "Laboratório de Tecido Conjuntivo / Instituto de Bioquímica Médica / Universidade Federal de Juiz de Fora"
$ws.Range("N23").Copy() | Out-Null
$ws.Range("O23").PasteSpecial(-4122) | Out-Null
$ws.Range("O23").Value = "Departamento de Glicobiologia  / Bioquímica Médica / Universidade Federal do Rio de Janeiro"
$ws.Range("P23").Copy() | Out-Null
$ws.Range("Q23").PasteSpecial(-4122) | Out-Null
$ws.Range("Q23").Value = "Universidade Federal do Rio de Janeiro"
# Row 24
$ws.Range("D24").Copy() | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").Value = "Department of Chemical and Biological Physics, Weizmann Institute of Science, Israel"
$ws.Range("F24").Copy() | Out-Null
$ws.Range("G24").PasteSpecial(-4122) | Out-Null
$ws.Range("G24").Value = "Centre National de La Recherche Scientifique, France"
$ws.Range("H24").Copy() | Out-Null
$ws.Range("I24").PasteSpecial(-4122) | Out-Null
$ws.Range("I24").Value = "Embrapa Instrumentation, Brazil"
$ws.Range("J24").Copy() | Out-Null
$ws.Range("K24").PasteSpecial(-4122) | Out-Null
$ws.Range("K24").Value = "São Carlos Institute of Physics, University of São Paulo, Brazil"
# Row 25
$ws.Range("D25").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = "Centre national de la recherche scientifique (CNRS)"
# Row 26
$ws.Range("D26").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "CEISAM / Faculté des Sciences et Techniques / Université de Nantes"
# Row 27
$ws.Range("D27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = "University of Toronto"
# Row 28
$ws.Range("D28").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "LowField"
# Row 29
$ws.Range("D29").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = "Universidade Federal de São Carlos"
$ws.Range("F29").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value = "Departamento de Química / Universidade Federal de São Carlos"
# Row 30
$ws.Range("D30").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = "NMR  / CIC energigune / CIC energigune"
# Row 31
$ws.Range("D31").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = "Dept Chem and BioM Engr / University of California"
# Row 32
$ws.Range("D32").Copy() | Out-Null
$ws.Range("E32").PasteSpecial(-4122) | Out-Null
$ws.Range("E32").Value = "University of Leipzig"
# Row 33
$ws.Range("D33").Copy() | Out-Null
$ws.Range("E33").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = "Universitat Autònoma de Barcelona"
# Row 34
$ws.Range("D34").Copy() | Out-Null
$ws.Range("E34").PasteSpecial(-4122) | Out-Null
$ws.Range("E34").Value = "NMR Research Scientist"
# Row 35
$ws.Range("D35").Copy() | Out-Null
$ws.Range("E35").PasteSpecial(-4122) | Out-Null
$ws.Range("E35").Value = "The University of Manchester"
# Row 36
$ws.Range("D36").Copy() | Out-Null
$ws.Range("E36").PasteSpecial(-4122) | Out-Null
$ws.Range("E36").Value = "Universidade Federal de São Carlos"
# Row 37
$ws.Range("D37").Copy() | Out-Null
$ws.Range("E37").PasteSpecial(-4122) | Out-Null
$ws.Range("E37").Value = "University of Toronto"
# Row 38
$ws.Range("D38").Copy() | Out-Null
$ws.Range("E38").PasteSpecial(-4122) | Out-Null
$ws.Range("E38").Value = "London, UK / Cryogenic Ltd"
$ws.Range("F38").Copy() | Out-Null
$ws.Range("G38").PasteSpecial(-4122) | Out-Null
$ws.Range("G38").Value = "Cryogenic Ltd"
$ws.Range("H38").Copy() | Out-Null
$ws.Range("I38").PasteSpecial(-4122) | Out-Null
$ws.Range("I38").Value = "Cryogenic Ltd"
# Row 39
$ws.Range("D39").Copy() | Out-Null
$ws.Range("E39").PasteSpecial(-4122) | Out-Null
$ws.Range("E39").Value = "Department of Chemistry, Materials and Chemical Engineering “Giulio Natta” / Politecnico di Milano"
# Row 40
$ws.Range("D40").Copy() | Out-Null
$ws.Range("E40").PasteSpecial(-4122) | Out-Null
$ws.Range("E40").Value = "Chimica / Chimica / Universidad Nacional de Rosario"
# Row 41
$ws.Range("D41").Copy() | Out-Null
$ws.Range("E41").PasteSpecial(-4122) | Out-Null
$ws.Range("E41").Value = "Chemistry / Arts & Sciences / Washington University in St. Louis"
$ws.Range("F41").Copy() | Out-Null
$ws.Range("G41").PasteSpecial(-4122) | Out-Null
$ws.Range("G41").Value = "Chemistry / Arts & Sciences / Washington University in St. Louis"
# Row 42
$ws.Range("D42").Copy() | Out-Null
$ws.Range("E42").PasteSpecial(-4122) | Out-Null
$ws.Range("E42").Value = "Universidade Federal do Espírito Santo, Brazil"
# Row 43
$ws.Range("D43").Copy() | Out-Null
$ws.Range("E43").PasteSpecial(-4122) | Out-Null
$ws.Range("E43").Value = "The University of Akron"
# Row 44
$ws.Range("D44").Copy() | Out-Null
$ws.Range("E44").PasteSpecial(-4122) | Out-Null
$ws.Range("E44").Value = "Iowa State University"
# Row 45
$ws.Range("D45").Copy() | Out-Null
$ws.Range("E45").PasteSpecial(-4122) | Out-Null
$ws.Range("E45").Value = "Niumag"
# Row 46
$ws.Range("D46").Copy() | Out-Null
$ws.Range("E46").PasteSpecial(-4122) | Out-Null
$ws.Range("E46").Value = "Nanalysis Corp."
# Row 47
$ws.Range("D47").Copy() | Out-Null
$ws.Range("E47").PasteSpecial(-4122) | Out-Null
$ws.Range("E47").Value = "Bruker"
# Row 48
$ws.Range("D48").Copy() | Out-Null
$ws.Range("E48").PasteSpecial(-4122) | Out-Null
$ws.Range("E48").Value = "University of Copenhagen"
$ws.Range("F48").Copy() | Out-Null
$ws.Range("G48").PasteSpecial(-4122) | Out-Null
$ws.Range("G48").Value = "University of Copenhagen"
$ws.Range("H48").Copy() | Out-Null
$ws.Range("I48").PasteSpecial(-4122) | Out-Null
$ws.Range("I48").Value = "University of Copenhagen"
# Row 49
$ws.Range("D49").Copy() | Out-Null
$ws.Range("E49").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Value = "Hélio para Ressonancia Magnética / Hélio para Ressonancia Magnética / Air Products do Brasil Ltda"
# Row 50
$ws.Range("D50").Copy() | Out-Null
$ws.Range("E50").PasteSpecial(-4122) | Out-Null
$ws.Range("E50").Value = "Oxford Instruments"
# Row 51
$ws.Range("D51").Copy() | Out-Null
$ws.Range("E51").PasteSpecial(-4122) | Out-Null
$ws.Range("E51").Value = "Magritek, Inc"
# Row 52
$ws.Range("D52").Copy() | Out-Null
$ws.Range("E52").PasteSpecial(-4122) | Out-Null
$ws.Range("E52").Value = "Universidade Federal de Ouro Preto"
# Row 53
$ws.Range("D53").Copy() | Out-Null
$ws.Range("E53").PasteSpecial(-4122) | Out-Null
$ws.Range("E53").Value = "INRAE Bordeaux, France"
# Row 54
$ws.Range("D54").Copy() | Out-Null
$ws.Range("E54").PasteSpecial(-4122) | Out-Null
$ws.Range("E54").Value = "Universidade de São Paulo"
$ws.Range("F54").Copy() | Out-Null
$ws.Range("G54").PasteSpecial(-4122) | Out-Null
$ws.Range("G54").Value = "Universidade de São Paulo"
$ws.Range("H54").Copy() | Out-Null
$ws.Range("I54").PasteSpecial(-4122) | Out-Null
$ws.Range("I54").Value = "Universidade Federal de São Carlos"
$ws.Range("J54").Copy() | Out-Null
$ws.Range("K54").PasteSpecial(-4122) | Out-Null
$ws.Range("K54").Value = "Universidade Federal de São Carlos"
$ws.Range("L54").Copy() | Out-Null
$ws.Range("M54").PasteSpecial(-4122) | Out-Null
$ws.Range("M54").Value = "Universidade Federal de São Carlos"
# Row 55
$ws.Range("D55").Copy() | Out-Null
$ws.Range("E55").PasteSpecial(-4122) | Out-Null
$ws.Range("E55").Value = "Laboratorio de Relaxometría y Técnicas Especiales (LaRTE). / Instituto de Física Enrique Gaviola, CONICET. Córdoba, Argentina. / FaMAF - Universidad Nacional de Córdoba. Córdoba, Argentina."
# Row 56
$ws.Range("D56").Copy() | Out-Null
$ws.Range("E56").PasteSpecial(-4122) | Out-Null
$ws.Range("E56").Value = "Chemistry / Scarborough Campus / University of Toronto"
$ws.Range("F56").Copy() | Out-Null
$ws.Range("G56").PasteSpecial(-4122) | Out-Null
$ws.Range("G56").Value = "University of Toronto"
$ws.Range("H56").Copy() | Out-Null
$ws.Range("I56").PasteSpecial(-4122) | Out-Null
$ws.Range("I56").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("J56").Copy() | Out-Null
$ws.Range("K56").PasteSpecial(-4122) | Out-Null
$ws.Range("K56").Value = "University of Toronto"
$ws.Range("L56").Copy() | Out-Null
$ws.Range("M56").PasteSpecial(-4122) | Out-Null
$ws.Range("M56").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("N56").Copy() | Out-Null
$ws.Range("O56").PasteSpecial(-4122) | Out-Null
$ws.Range("O56").Value = "University of Toronto"
$ws.Range("P56").Copy() | Out-Null
$ws.Range("Q56").PasteSpecial(-4122) | Out-Null
$ws.Range("Q56").Value = "Bruker"
$ws.Range("R56").Copy() | Out-Null
$ws.Range("S56").PasteSpecial(-4122) | Out-Null
$ws.Range("S56").Value = "Biospin GmbH / BRUKER"
$ws.Range("T56").Copy() | Out-Null
$ws.Range("U56").PasteSpecial(-4122) | Out-Null
$ws.Range("U56").Value = "Biospin AG / BRUKER"
$ws.Range("V56").Copy() | Out-Null
$ws.Range("W56").PasteSpecial(-4122) | Out-Null
$ws.Range("W56").Value = "Biospin AG / BRUKER"
$ws.Range("X56").Copy() | Out-Null
$ws.Range("Y56").PasteSpecial(-4122) | Out-Null
$ws.Range("Y56").Value = "Biospin AG / BRUKER"
$ws.Range("Z56").Copy() | Out-Null
$ws.Range("AA56").PasteSpecial(-4122) | Out-Null
$ws.Range("AA56").Value = "Biospin AG / BRUKER"
$ws.Range("AB56").Copy() | Out-Null
$ws.Range("AC56").PasteSpecial(-4122) | Out-Null
$ws.Range("AC56").Value = "Biospin AG / BRUKER"
$ws.Range("AD56").Copy() | Out-Null
$ws.Range("AE56").PasteSpecial(-4122) | Out-Null
$ws.Range("AE56").Value = "Biospin AG / BRUKER"
$ws.Range("AF56").Copy() | Out-Null
$ws.Range("AG56").PasteSpecial(-4122) | Out-Null
$ws.Range("AG56").Value = "Canada Ltd / BRUKER"
$ws.Range("AH56").Copy() | Out-Null
$ws.Range("AI56").PasteSpecial(-4122) | Out-Null
$ws.Range("AI56").Value = "Canada Ltd / BRUKER"
$ws.Range("AJ56").Copy() | Out-Null
$ws.Range("AK56").PasteSpecial(-4122) | Out-Null
$ws.Range("AK56").Value = "University of Toronto"
$ws.Range("AL56").Copy() | Out-Null
$ws.Range("AM56").PasteSpecial(-4122) | Out-Null
$ws.Range("AM56").Value = "University of Toronto"
# Row 57
$ws.Range("D57").Copy() | Out-Null
$ws.Range("E57").PasteSpecial(-4122) | Out-Null
$ws.Range("E57").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("F57").Copy() | Out-Null
$ws.Range("G57").PasteSpecial(-4122) | Out-Null
$ws.Range("G57").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("H57").Copy() | Out-Null
$ws.Range("I57").PasteSpecial(-4122) | Out-Null
$ws.Range("I57").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("J57").Copy() | Out-Null
$ws.Range("K57").PasteSpecial(-4122) | Out-Null
$ws.Range("K57").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("L57").Copy() | Out-Null
$ws.Range("M57").PasteSpecial(-4122) | Out-Null
$ws.Range("M57").Value = "Biospin AG / Bruker"
$ws.Range("N57").Copy() | Out-Null
$ws.Range("O57").PasteSpecial(-4122) | Out-Null
$ws.Range("O57").Value = "Biospin AG / Bruker"
$ws.Range("P57").Copy() | Out-Null
$ws.Range("Q57").PasteSpecial(-4122) | Out-Null
$ws.Range("Q57").Value = "Biospin AG / Bruker"
$ws.Range("R57").Copy() | Out-Null
$ws.Range("S57").PasteSpecial(-4122) | Out-Null
$ws.Range("S57").Value = "Bruker Biospin GmbH / Bruker"
$ws.Range("T57").Copy() | Out-Null
$ws.Range("U57").PasteSpecial(-4122) | Out-Null
$ws.Range("U57").Value = "Bruker Biospin GmbH / Bruker"
$ws.Range("V57").Copy() | Out-Null
$ws.Range("W57").PasteSpecial(-4122) | Out-Null
$ws.Range("W57").Value = "Bruker BioSpin"
$ws.Range("X57").Copy() | Out-Null
$ws.Range("Y57").PasteSpecial(-4122) | Out-Null
$ws.Range("Y57").Value = "Bruker"
$ws.Range("Z57").Copy() | Out-Null
$ws.Range("AA57").PasteSpecial(-4122) | Out-Null
$ws.Range("AA57").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
$ws.Range("AB57").Copy() | Out-Null
$ws.Range("AC57").PasteSpecial(-4122) | Out-Null
$ws.Range("AC57").Value = "Department of Physical and Environmental Science / Scarborough Campus / University of Toronto"
# Row 58
$ws.Range("D58").Copy() | Out-Null
$ws.Range("E58").PasteSpecial(-4122) | Out-Null
$ws.Range("E58").Value = "Universidade Federal do Estado do Rio de Janeiro"
$ws.Range("F58").Copy() | Out-Null
$ws.Range("G58").PasteSpecial(-4122) | Out-Null
$ws.Range("G58").Value = "Universidade Federal do Rio de Janeiro"
$ws.Range("H58").Copy() | Out-Null
$ws.Range("I58").PasteSpecial(-4122) | Out-Null
$ws.Range("I58").Value = "Universidade Federal do Rio de Janeiro"
$ws.Range("J58").Copy() | Out-Null
$ws.Range("K58").PasteSpecial(-4122) | Out-Null
$ws.Range("K58").Value = "Universidade Federal do Rio de Janeiro"
$ws.Range("L58").Copy() | Out-Null
$ws.Range("M58").PasteSpecial(-4122) | Out-Null
$ws.Range("M58").Value = "Universidade Federal do Rio de Janeiro"
$ws.Range("N58").Copy() | Out-Null
$ws.Range("O58").PasteSpecial(-4122) | Out-Null
$ws.Range("O58").Value = "Universidade Federal do Estado do Rio de Janeiro"
# Row 59
$ws.Range("D59").Copy() | Out-Null
$ws.Range("E59").PasteSpecial(-4122) | Out-Null
$ws.Range("E59").Value = "Universidade Federal do Rio de Janeiro"
# Row 60
$ws.Range("D60").Copy() | Out-Null
$ws.Range("E60").PasteSpecial(-4122) | Out-Null
$ws.Range("E60").Value = "NMR Applicaitons / Bruker Biospin"
# Row 61
$ws.Range("D61").Copy() | Out-Null
$ws.Range("E61").PasteSpecial(-4122) | Out-Null
$ws.Range("E61").Value = "Department of Pharmaceutical Sciences / Faculty of Life Sciences / University of Vienna"
# Row 62
$ws.Range("D62").Copy() | Out-Null
$ws.Range("E62").PasteSpecial(-4122) | Out-Null
$ws.Range("E62").Value = "Alegre Science Inc"
# Row 63
$ws.Range("D63").Copy() | Out-Null
$ws.Range("E63").PasteSpecial(-4122) | Out-Null
$ws.Range("E63").Value = "Department of Chemistry, Materials and Chemical Engineering “Giulio Natta” / Politecnico di Milano"
$ws.Range("F63").Copy() | Out-Null
$ws.Range("G63").PasteSpecial(-4122) | Out-Null
$ws.Range("G63").Value = "Department of Chemistry, Materials and Chemical Engineering “Giulio Natta” / Politecnico di Milano"
$ws.Range("H63").Copy() | Out-Null
$ws.Range("I63").PasteSpecial(-4122) | Out-Null
$ws.Range("I63").Value = "Department for Sustainability (SSPT) / Italian National Agency for New Technologies, Energy and Sustainable Economic Development / ENEA"
$ws.Range("J63").Copy() | Out-Null
$ws.Range("K63").PasteSpecial(-4122) | Out-Null
$ws.Range("K63").Value = "School of Chemical Engineering/ University of Campinas/ Campinas/ Brazil"
$ws.Range("L63").Copy() | Out-Null
$ws.Range("M63").PasteSpecial(-4122) | Out-Null
$ws.Range("M63").Value = "Department of Chemistry, Materials and Chemical Engineering “Giulio Natta” / Politecnico di Milano"
$excel.CutCopyMode = 0
